$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the loading details column header
$ws.Range("P8").Value = "24V Rail(A)"

# Update the slot-card device label text to reflect new amperage rating
$ws.Range("G8").Value = "PCH800 5.0A,PNI800"
$ws.Range("H8").Value = "PCH800 5.0A-1,PNI800-2"

# Update the visible/active selection to match latest edit position
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("L8").Select()
